$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.898.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.035.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.34"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.335.10"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.54"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.013.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.827.33"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.82"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.90"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0606"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.42"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.524.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.04"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.05"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0916"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.10"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.223.88"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.82%  "
